# Update the answer values in the "two-digit x two-digit" multiplication
# table. Each populated row holds five "a×b=c" answers; only the answer
# text itself changes, so each cell is addressed explicitly by its
# (row, column) position in Table 1 and its Range.Text is replaced,
# leaving every other run/paragraph property untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Old="66×73=4818"; New="11×96=1056"},
    @{Row=1;  Col=2; Old="65×28=1820"; New="84×66=5544"},
    @{Row=1;  Col=3; Old="66×29=1914"; New="76×36=2736"},
    @{Row=1;  Col=4; Old="97×76=7372"; New="51×41=2091"},
    @{Row=1;  Col=5; Old="89×40=3560"; New="17×92=1564"},

    @{Row=5;  Col=1; Old="52×80=4160"; New="75×59=4425"},
    @{Row=5;  Col=2; Old="64×68=4352"; New="15×73=1095"},
    @{Row=5;  Col=3; Old="54×25=1350"; New="58×45=2610"},
    @{Row=5;  Col=4; Old="53×46=2438"; New="49×53=2597"},
    @{Row=5;  Col=5; Old="35×11=385";  New="24×97=2328"},

    @{Row=10; Col=1; Old="52×62=3224"; New="44×81=3564"},
    @{Row=10; Col=2; Old="39×60=2340"; New="46×48=2208"},
    @{Row=10; Col=3; Old="37×37=1369"; New="14×19=266"},
    @{Row=10; Col=4; Old="78×91=7098"; New="89×38=3382"},
    @{Row=10; Col=5; Old="76×36=2736"; New="15×68=1020"},

    @{Row=15; Col=1; Old="89×18=1602"; New="39×90=3510"},
    @{Row=15; Col=2; Old="47×97=4559"; New="40×60=2400"},
    @{Row=15; Col=3; Old="78×75=5850"; New="43×32=1376"},
    @{Row=15; Col=4; Old="46×83=3818"; New="25×67=1675"},
    @{Row=15; Col=5; Old="41×17=697";  New="99×97=9603"},

    @{Row=20; Col=1; Old="82×30=2460"; New="20×79=1580"},
    @{Row=20; Col=2; Old="33×93=3069"; New="49×93=4557"},
    @{Row=20; Col=3; Old="65×32=2080"; New="53×17=901"},
    @{Row=20; Col=4; Old="36×95=3420"; New="72×18=1296"},
    @{Row=20; Col=5; Old="89×73=6497"; New="34×67=2278"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $range = $cell.Range
    # Trim the cell-end marker off before comparing so we match just the
    # visible answer text (the Range includes a trailing \r\a pair).
    $current = $range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $u.Old) {
        throw "Unexpected text at row $($u.Row) col $($u.Col): [$current] (expected [$($u.Old)])"
    }
    $range.Text = $u.New
}

Write-Output "Updated $($updates.Count) answer cells"
